# Update countries & provincias Spain
# Applies the refreshed COVID-19 snapshot values to paises.xlsx:
#   - timestamp string in A1
#   - three pairs of countries that swapped ranking order (name + numbers)
#   - numeric refreshes for a handful of other countries

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $vals) {
    if ($country -ne $null) {
        $ws.Cells.Item($row, 1).Value = $country
    }
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
    $ws.Cells.Item($row, 8).Value = $vals[6]
}

# --- Row 1: "Datos actualizados ..." timestamp ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 24 de Agosto de 2020 a las 04:29"

# --- Rows 29 (Bolivia), 72 (Australia), 84 (Paraguay), 143 (Nueva Zelanda): straight numeric refresh ---
Set-Row 29  $null @(109149, 722, 45396, 59244, 0, 67, 4509)
Set-Row 72  $null @(24915, 103, 19233, 5165, 0, 15, 517)
Set-Row 84  $null @(13233, 0, 7417, 5611, 0, 0, 205)
Set-Row 143 $null @(1683, 9, 1538, 123, 0, 0, 22)

# --- Rows 78/79: Corea del Sur overtakes Costa de Marfil ---
Set-Row 78 "Corea del Sur"   @(17665, 266, 14219, 3137, 0, 0, 309)
Set-Row 79 "Costa de Marfil" @(17471, 0, 15301, 2057, 0, 0, 113)

# --- Rows 147/148/149: Jamaica overtakes Uruguay and Republica de Chipre ---
Set-Row 147 "Jamaica"             @(1529, 116, 819, 694, 0, 0, 16)
Set-Row 148 "Uruguay"             @(1527, 0, 1276, 209, 0, 0, 42)
Set-Row 149 "Republica de Chipre" @(1421, 0, 878, 523, 0, 0, 20)

# --- Rows 202/203: Timor Oriental overtakes Santa Lucia (values identical, only names swap) ---
Set-Row 202 "Timor Oriental" @(26, 0, 25, 1, 0, 0, 0)
Set-Row 203 "Santa Lucia"    @(26, 0, 25, 1, 0, 0, 0)
